# Update market/profit data values in several Leve sheets
# (scheduled market-data refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 270.2
$ws.Range("I12").Value = 216.33333
$ws.Range("K12").Value = 216.33333
$ws.Range("M12").Value = -46.33332999999999
$ws.Range("H19").Value = 1085.7693
$ws.Range("I19").Value = 1359.9
$ws.Range("J19").Value = 172
$ws.Range("K19").Value = 1359.9
$ws.Range("L19").Value = 172
$ws.Range("M19").Value = -1184.9
$ws.Range("N19").Value = -522
$ws.Range("H47").Value = 5000
$ws.Range("I47").Value = 5000
$ws.Range("K47").Value = 5000
$ws.Range("M47").Value = -4028
$ws.Range("H58").Value = 347.5
$ws.Range("J58").Value = 1000
$ws.Range("L58").Value = 3000
$ws.Range("N58").Value = -3300
$ws.Range("H62").Value = 2456.5715
$ws.Range("I62").Value = 1999.3334
$ws.Range("J62").Value = 5200
$ws.Range("K62").Value = 1999.3334
$ws.Range("L62").Value = 5200
$ws.Range("M62").Value = -1375.3334
$ws.Range("N62").Value = -6448
$ws.Range("H65").Value = 2456.5715
$ws.Range("I65").Value = 1999.3334
$ws.Range("J65").Value = 5200
$ws.Range("K65").Value = 9996.666999999999
$ws.Range("L65").Value = 26000
$ws.Range("M65").Value = -6876.666999999999
$ws.Range("N65").Value = -32240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 78
$ws.Range("I5").Value = 46.8
$ws.Range("J5").Value = 156
$ws.Range("K5").Value = 46.8
$ws.Range("L5").Value = 156
$ws.Range("M5").Value = 65.2
$ws.Range("N5").Value = -380
$ws.Range("H102").Value = 771.8
$ws.Range("I102").Value = 640
$ws.Range("K102").Value = 640
$ws.Range("M102").Value = 982

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 78
$ws.Range("I4").Value = 46.8
$ws.Range("J4").Value = 156
$ws.Range("K4").Value = 46.8
$ws.Range("L4").Value = 156
$ws.Range("M4").Value = 68.2
$ws.Range("N4").Value = -386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 12999
$ws.Range("I58").Value = 8000
$ws.Range("K58").Value = 8000
$ws.Range("M58").Value = -7797
$ws.Range("H59").Value = 27658.334
$ws.Range("J59").Value = 44980
$ws.Range("L59").Value = 44980
$ws.Range("N59").Value = -47270
$ws.Range("H107").Value = 1500
$ws.Range("H136").Value = 12999
$ws.Range("I136").Value = 8000
$ws.Range("K136").Value = 24000
$ws.Range("M136").Value = -21450
$ws.Range("H141").Value = 424800.1
$ws.Range("J141").Value = 424800.1
$ws.Range("L141").Value = 424800.1
$ws.Range("N141").Value = -435160.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 422105.75
$ws.Range("I4").Value = 1004010
$ws.Range("J4").Value = 6459.857
$ws.Range("K4").Value = 3012030
$ws.Range("L4").Value = 19379.571
$ws.Range("M4").Value = -3011918
$ws.Range("N4").Value = -19603.571
$ws.Range("H12").Value = 44.384617
$ws.Range("I12").Value = 56.5
$ws.Range("J12").Value = 39
$ws.Range("K12").Value = 169.5
$ws.Range("L12").Value = 117
$ws.Range("M12").Value = 3.5
$ws.Range("N12").Value = -463
$ws.Range("H40").Value = 179.8
$ws.Range("J40").Value = 74
$ws.Range("L40").Value = 296
$ws.Range("N40").Value = -434
$ws.Range("H61").Value = 999
$ws.Range("J61").Value = 999
$ws.Range("L61").Value = 2997
$ws.Range("N61").Value = -3427
$ws.Range("H114").Value = 1530
$ws.Range("I114").Value = 1530
$ws.Range("K114").Value = 4590
$ws.Range("M114").Value = -1336
$ws.Range("H117").Value = 3165
$ws.Range("J117").Value = 3165
$ws.Range("L117").Value = 9495
$ws.Range("N117").Value = -16379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10003
$ws.Range("I10").Value = 10003
$ws.Range("K10").Value = 10003
$ws.Range("M10").Value = -9834
$ws.Range("H12").Value = 15004
$ws.Range("J12").Value = 15004
$ws.Range("L12").Value = 15004
$ws.Range("N12").Value = -15284
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H122").Value = 950
$ws.Range("I122").Value = 940
$ws.Range("K122").Value = 2820
$ws.Range("M122").Value = -370

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 933.1111
$ws.Range("I22").Value = 833.1667
$ws.Range("K22").Value = 833.1667
$ws.Range("M22").Value = -538.1667
$ws.Range("H27").Value = 933.1111
$ws.Range("I27").Value = 833.1667
$ws.Range("K27").Value = 833.1667
$ws.Range("M27").Value = -726.1667
$ws.Range("H46").Value = 8493
$ws.Range("I46").Value = 8708.5
$ws.Range("K46").Value = 8708.5
$ws.Range("M46").Value = -8520.5
$ws.Range("H122").Value = 3196
$ws.Range("I122").Value = 2995
$ws.Range("K122").Value = 8985
$ws.Range("M122").Value = -6535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H46").Value = 99999
$ws.Range("J46").Value = 99999
$ws.Range("L46").Value = 99999
$ws.Range("N46").Value = -100461
$ws.Range("H62").Value = 5166.6665
$ws.Range("J62").Value = 6500
$ws.Range("L62").Value = 6500
$ws.Range("N62").Value = -7748
$ws.Range("H65").Value = 5166.6665
$ws.Range("J65").Value = 6500
$ws.Range("L65").Value = 32500
$ws.Range("N65").Value = -38740
$ws.Range("H134").Value = 99999
$ws.Range("J134").Value = 99999
$ws.Range("L134").Value = 299997
$ws.Range("N134").Value = -305067
